# Update Data by bot, scripted by HH
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / inline-string fields (row 2) ---
# J2 looks numeric ("001"); prefix with an apostrophe so Excel stores it
# as literal text instead of coercing it to the number 1.
$ws.Range("J2").Value = "'001"

$ws.Range("M2").Value = "2020-12-22 00:00:00"
$ws.Range("N2").Value = "2017-12-31 00:00:00"

# --- Numeric fields (row 2) ---
$ws.Range("O2").Value = 39055812.29
$ws.Range("P2").Value = 169579892.56
$ws.Range("Q2").Value = 124337606.5
$ws.Range("R2").Value = 23.4852263003
$ws.Range("S2").Value = 107477807.38
$ws.Range("T2").Value = 107477807.38
$ws.Range("U2").Value = 24.996617495
$ws.Range("V2").Value = 1486637.68
$ws.Range("W2").Value = 6595069.59
$ws.Range("X2").Value = -393563.88
$ws.Range("Y2").Value = 45954759.22
$ws.Range("Z2").Value = 45954733.28
$ws.Range("AA2").Value = 6898920.99
$ws.Range("AG2").Value = 825137.16
$ws.Range("AP2").Value = 22.7768406497
$ws.Range("AQ2").Value = 25.970437614312
$ws.Range("AR2").Value = 22.737074754653
$ws.Range("AS2").Value = 37791846.04
$ws.Range("AT2").Value = 23.19679323622
